# REM22062023_CARRO.xlsx - switch the xlsx reader from ClosedXML to MiniExcel.
# Re-saving through the new library normalizes the header row casing and
# nudges a couple of cosmetic sheet settings; reproduce those observable
# changes here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers were lower-case ("product"/"customer"/"brand"/"value"); capitalize them.
$ws.Range("A1").Value = "Product"
$ws.Range("B1").Value = "Customer"
$ws.Range("C1").Value = "Brand"
$ws.Range("D1").Value = "Value"

# Column C (Brand) gets an explicit best-fit-sized width.
$ws.Columns.Item(3).ColumnWidth = 6.166666666666667

# H5 picks up an underline font (the cell that ends up selected below).
$ws.Range("H5").Font.Underline = $true

# Last active selection moves to H5.
$ws.Range("H5").Select()

# Sheet is set up for A4 portrait printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
